# Update the "Förändrad" (Changed) date column (C) for all data rows
# from serial date 45170 (2023-09-01) to 45174 (2023-09-05).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

$newValue = 45174

for ($row = 2; $row -le 41; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45170) {
        $cell.Value = $newValue
    }
}
